# Apply updated symbol list values (price & 1h volume change) to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values are written as text (matching the
# original inline-string cells) rather than being auto-converted to numbers/percentages
# by Excel's automatic type detection.
$updates = @{
    "D2" = '298.56'
    "E2" = '-1.84%'
    "D3" = '31.24'
    "E3" = '-1.25%'
    "D4" = '5.119'
    "E4" = '-0.91%'
    "D5" = '0.07949'
    "E5" = '6.07%'
    "D6" = '2.230'
    "E6" = '-6.75%'
    "D7" = '7.777'
    "E7" = '-2.97%'
    "D8" = '3.861'
    "E8" = '-0.14%'
    "D9" = '0.9220'
    "E9" = '0.74%'
    "D10" = '0.1727'
    "E10" = '-0.38%'
    "D11" = '0.07488'
    "E11" = '-2.77%'
    "D12" = '0.09237'
    "E12" = '13.11%'
    "D13" = '0.03064'
    "E13" = '0.83%'
    "D14" = '0.1003'
    "E14" = '0.97%'
    "D15" = '0.001513'
    "E15" = '0.51%'
    "D16" = '0.006026'
    "E16" = '-1.06%'
    "E17" = '-0.63%'
    "E18" = '1.31%'
    "E19" = '0.25%'
    "D20" = '0.1334'
    "E20" = '0.34%'
    "D21" = '3.913'
    "E21" = '-15.91%'
    "D22" = '0.1699'
    "E22" = '8.61%'
    "D23" = '0.04620'
    "E23" = '0.56%'
    "D24" = '0.001245'
    "E24" = '-1.23%'
    "D25" = '0.004480'
    "E25" = '-1.24%'
    "E26" = '-7.57%'
    "D27" = '0.0003396'
    "E27" = '24.06%'
    "D39" = '0.01756'
    "E39" = '0.14%'
    "D40" = '0.04601'
    "E40" = '1.34%'
    "E41" = '-5.91%'
    "D42" = '0.1361'
    "E42" = '-0.20%'
    "D43" = '0.002189'
    "E43" = '1.06%'
    "E44" = '-7.60%'
    "D45" = '0.00006314'
    "E45" = '-3.51%'
    "D46" = '0.00000000750'
    "E46" = '-0.07%'
    "D47" = '0.007975'
    "E47" = '-19.32%'
    "D48" = '0.7465'
    "E48" = '-9.03%'
    "D49" = '0.00002099'
    "E49" = '-0.07%'
    "D50" = '0.0001999'
    "E50" = '0.00%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
